# Update "想去人数" (F column) counts and one "最低票价" (G15) value
# across the "展览" and "全部类型" sheets, mirroring the freshly
# regenerated gh-pages data output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 161
    $ws.Range("F3").Value = 7280
    $ws.Range("F4").Value = 5482
    $ws.Range("F6").Value = 171
    $ws.Range("F7").Value = 34
    $ws.Range("F8").Value = 44
    $ws.Range("F9").Value = 107
    $ws.Range("F10").Value = 84
    $ws.Range("F11").Value = 104
    $ws.Range("F12").Value = 203
    $ws.Range("F13").Value = 35
    $ws.Range("F14").Value = 647

    # Row 15 became sellable: a ticket price replaces "不可售"
    $ws.Range("F15").Value = 275
    $ws.Range("G15").Value = 65

    $ws.Range("F16").Value = 52
    $ws.Range("F17").Value = 6
    $ws.Range("F19").Value = 33
}
